# Natmi following Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values per row (rows 2-5), columns E,G,H,K,M,N,O,P,Q,R,S,T
$updates = @(
    @{ Row = 2; E = 3; G = 0.3430696666666667; H = 1.029209; K = 3; M = 75.59011833333334; N = 226.770355; O = 0.6588374259037486; P = 0.6588374259037486; Q = 25.93267669991056; R = 233.394090299195; S = 0.6588374259037486; T = 0.6588374259037486 },
    @{ Row = 3; E = 3; G = 0.3430696666666667; H = 1.029209; K = 3; M = 12.15310033333333; N = 36.459301; O = 0.1059254505338229; P = 0.1059254505338229; Q = 4.169360080323222; R = 37.524240722909; S = 0.1059254505338229; T = 0.1059254505338229 },
    @{ Row = 4; E = 3; G = 0.3430696666666667; H = 1.029209; K = 3; M = 0.1727356666666667; N = 0.5182070000000001; O = 0.001505550255743542; P = 0.001505550255743542; Q = 0.05926036758477779; R = 0.5333433082630001; S = 0.001505550255743542; T = 0.001505550255743542 },
    @{ Row = 5; E = 3; G = 0.3430696666666667; H = 1.029209; K = 3; M = 26.81662666666667; N = 80.44988000000001; O = 0.233731573306685; P = 0.233731573306685; Q = 9.199971171657777; R = 82.79974054492001; S = 0.233731573306685; T = 0.233731573306685 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("E$r").Value = $u.E
    $ws.Range("G$r").Value = $u.G
    $ws.Range("H$r").Value = $u.H
    $ws.Range("K$r").Value = $u.K
    $ws.Range("M$r").Value = $u.M
    $ws.Range("N$r").Value = $u.N
    $ws.Range("O$r").Value = $u.O
    $ws.Range("P$r").Value = $u.P
    $ws.Range("Q$r").Value = $u.Q
    $ws.Range("R$r").Value = $u.R
    $ws.Range("S$r").Value = $u.S
    $ws.Range("T$r").Value = $u.T
}
